$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text so Excel does not
# auto-convert decimal-looking strings (e.g. "1.004") into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.912.91"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.691.71"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +1.00%  "
$ws.Range("D5").Value = "315.49"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "0.3952"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").Value = "0.3997"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").Value = "1.445"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").Value = "52.36"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "0.08727"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "25.51"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").Value = "7.389"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "0.00001338"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "7.871"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").Value = "1.675.44"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "94.82"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").Value = "0.07196"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "7.182"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "1.005"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "14.18"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "24.919.59"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").Value = "2.406"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "2.852"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("D27").Value = "23.07"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "6.050"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").Value = "162.46"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("D30").Value = "148.49"
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("B31").Value = "WEMIXTOKEN"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "2.631"
$ws.Range("E31").Value = "  +21.46%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "8.032"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").Value = "1.926.65"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "0.08517"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("D35").Value = "0.03113"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "1.035"
$ws.Range("E36").Value = "  -2.23%  "
$ws.Range("D37").Value = "7.024"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").Value = "0.2865"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "0.09700"
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "0.8090"
$ws.Range("E41").Value = "  -7.46%  "
$ws.Range("D42").Value = "13.89"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").Value = "1.472"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "16.95"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "0.7283"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").Value = "4.222"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "0.08937"
$ws.Range("E48").Value = "  +8.24%  "
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "1.006"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "139.23"
$ws.Range("E51").Value = "  -1.14%  "

# Restore the original (default) cell style on the Price column so
# no stray number-format / quote-prefix style lingers on the cells.
$priceRange.Style = "Normal"

